# IST price update 2025-12-21 13:39
#
# A fresh price-scrape column is inserted immediately to the right of the
# "SKU Name" column (new column B), pushing the previous columns B:N one
# slot to the right (C:O). The new column starts out as a duplicate of the
# scrape it is based on (the prior column B, which after the shift now
# lives in column C) and is then corrected with the timestamp for this run
# plus the one SKU whose price actually changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 26

# Remember column B's current width so the freshly inserted column keeps
# the same "21" width as every other price column instead of the default.
$colWidth = $ws.Columns("B:B").ColumnWidth()

# Shift B:N -> C:O, leaving a blank column B behind.
$ws.Columns("B:B").Insert()
$ws.Columns("B:B").ColumnWidth = $colWidth

# The new column starts as a copy of the run it was cloned from (now in C).
$ws.Range("B1:B$lastRow").Value = $ws.Range("C1:C$lastRow").Value()

# New scrape timestamp for this run's header.
$ws.Range("B1").Value = "2025-12-21 19:06"

# Row 19 ("...162 Count | Pack of 3...") price dropped in this run.
$ws.Range("B19").Value = 465
